$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: set J7 to "000010" (new shared string), then I7 to 0
$ws.Range("J7").Value = "000010"
$ws.Range("I7").Value = 0

# Row 8: set J8 to "000001" (new shared string)
$ws.Range("J8").Value = "000001"

# Row 7: add A7 value "B1" (new shared string, last one added)
$ws.Range("A7").Value = "B1"

# K10: 5 -> 4.5
$ws.Range("K10").Value = 4.5

# L11: 3 -> 2.5
$ws.Range("L11").Value = 2.5

# Update selection to I8
$ws.Range("I8").Select()
